$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'30.711.17"
$ws.Range("E2").Value = "'  +0.23%  "

# Row 3
$ws.Range("D3").Value = "'1.920.39"
$ws.Range("E3").Value = "'  +1.62%  "

# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "'  +0.16%  "

# Row 5
$ws.Range("D5").Value = "'240.67"
$ws.Range("E5").Value = "'  -1.85%  "

# Row 6
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "'  +0.12%  "

# Row 7
$ws.Range("D7").Value = "'0.4933"
$ws.Range("E7").Value = "'  +0.32%  "

# Row 8
$ws.Range("D8").Value = "'0.3005"
$ws.Range("E8").Value = "'  +1.52%  "

# Row 9
$ws.Range("D9").Value = "'0.06777"
$ws.Range("E9").Value = "'  -0.34%  "

# Row 10
$ws.Range("D10").Value = "'1.899.92"
$ws.Range("E10").Value = "'  +0.57%  "

# Row 11
$ws.Range("D11").Value = "'17.28"
$ws.Range("E11").Value = "'  +0.05%  "

# Row 12
$ws.Range("D12").Value = "'0.07338"
$ws.Range("E12").Value = "'  +1.27%  "

# Row 13
$ws.Range("D13").Value = "'5.220"
$ws.Range("E13").Value = "'  +3.46%  "

# Row 14
$ws.Range("D14").Value = "'88.77"
$ws.Range("E14").Value = "'  -2.76%  "

# Row 15
$ws.Range("D15").Value = "'0.6761"
$ws.Range("E15").Value = "'  -0.24%  "

# Row 16
$ws.Range("D16").Value = "'30.694.03"
$ws.Range("E16").Value = "'  +0.25%  "

# Row 17
$ws.Range("D17").Value = "'0.000007979"
$ws.Range("E17").Value = "'  +0.12%  "

# Row 18
$ws.Range("D18").Value = "'13.58"
$ws.Range("E18").Value = "'  +2.98%  "

# Row 19
$ws.Range("E19").Value = "'  +0.07%  "

# Row 20
$ws.Range("D20").Value = "'2.162.61"
$ws.Range("E20").Value = "'  +1.51%  "

# Row 21
$ws.Range("D21").Value = "'5.423"
$ws.Range("E21").Value = "'  +12.40%  "

# Row 22
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "'  +0.13%  "

# Row 23
$ws.Range("D23").Value = "'198.56"
$ws.Range("E23").Value = "'  +2.64%  "

# Row 24
$ws.Range("D24").Value = "'6.373"
$ws.Range("E24").Value = "'  +4.88%  "

# Row 25
$ws.Range("D25").Value = "'9.682"
$ws.Range("E25").Value = "'  +3.62%  "

# Row 26
$ws.Range("D26").Value = "'164.28"
$ws.Range("E26").Value = "'  +5.93%  "

# Row 27
$ws.Range("D27").Value = "'18.73"
$ws.Range("E27").Value = "'  -2.64%  "

# Row 28
$ws.Range("D28").Value = "'1.968"
$ws.Range("E28").Value = "'  +3.19%  "

# Row 29
$ws.Range("E29").Value = "'  +4.69%  "

# Row 30
$ws.Range("D30").Value = "'4.383"
$ws.Range("E30").Value = "'  +1.08%  "

# Row 31
$ws.Range("D31").Value = "'0.09189"
$ws.Range("E31").Value = "'  +1.12%  "

# Row 32
$ws.Range("D32").Value = "'4.078"
$ws.Range("E32").Value = "'  +1.45%  "

# Row 33
$ws.Range("D33").Value = "'0.05284"
$ws.Range("E33").Value = "'  +1.53%  "

# Row 34
$ws.Range("D34").Value = "'0.7455"
$ws.Range("E34").Value = "'  -2.17%  "

# Row 35
$ws.Range("D35").Value = "'1.121"
$ws.Range("E35").Value = "'  +0.89%  "

# Row 36
$ws.Range("D36").Value = "'2.710"
$ws.Range("E36").Value = "'  -2.28%  "

# Row 37
$ws.Range("D37").Value = "'0.01855"
$ws.Range("E37").Value = "'  +0.81%  "

# Row 38
$ws.Range("D38").Value = "'2.729"
$ws.Range("E38").Value = "'  +2.06%  "

# Row 39
$ws.Range("D39").Value = "'0.9284"
$ws.Range("E39").Value = "'  -0.74%  "

# Row 40
$ws.Range("D40").Value = "'2.086"
$ws.Range("E40").Value = "'  -3.04%  "

# Row 41
$ws.Range("D41").Value = "'0.4499"
$ws.Range("E41").Value = "'  +1.63%  "

# Row 42
$ws.Range("D42").Value = "'72.37"
$ws.Range("E42").Value = "'  +25.07%  "

# Row 43
$ws.Range("D43").Value = "'5.967"
$ws.Range("E43").Value = "'  +3.64%  "

# Row 44
$ws.Range("D44").Value = "'106.96"
$ws.Range("E44").Value = "'  +1.46%  "

# Row 45
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = "'1.003"
$ws.Range("E45").Value = "'  +0.30%  "

# Row 46
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = "'0.1406"
$ws.Range("E46").Value = "'  +4.30%  "

# Row 47
$ws.Range("D47").Value = "'7.686"
$ws.Range("E47").Value = "'  +1.10%  "

# Row 48
$ws.Range("D48").Value = "'9.063"
$ws.Range("E48").Value = "'  +4.21%  "

# Row 49
$ws.Range("D49").Value = "'35.36"
$ws.Range("E49").Value = "'  +5.31%  "

# Row 50
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = "'0.05900"
$ws.Range("E50").Value = "'  +0.67%  "

# Row 51
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = "'0.4059"
$ws.Range("E51").Value = "'  +3.28%  "
